# Weekly update: insert a new price record as row 112 for
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Zapallo", pushing the
# existing rows 112-139 down to 113-140.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 112 (shifts 112..139 -> 113..140,
# carries the D-column date style down automatically).
$ws.Rows.Item(112).Insert()

# Populate the new row 112 with this week's record.
$ws.Cells.Item(112, 1).Value  = 7
$ws.Cells.Item(112, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(112, 3).Value  = "Ñuble"
$ws.Cells.Item(112, 4).Value  = 44642
$ws.Cells.Item(112, 5).Value  = 16
$ws.Cells.Item(112, 6).Value  = 100112045
$ws.Cells.Item(112, 7).Value  = "Zapallo"
$ws.Cells.Item(112, 8).Value  = "Camote"
$ws.Cells.Item(112, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(112, 10).Value = 120
$ws.Cells.Item(112, 11).Value = 300
$ws.Cells.Item(112, 12).Value = 350
$ws.Cells.Item(112, 13).Value = 325
$ws.Cells.Item(112, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(112, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(112, 16).Value = 325
$ws.Cells.Item(112, 17).Value = 1
$ws.Cells.Item(112, 18).Value = "Hortaliza"
